$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = 0.68836793960366061
$ws.Range("AD1").Value = 0.96325046501457168
$ws.Range("B3").Value = 0.97861784562256326
$ws.Range("M3").Value = 0.78887886096553261
$ws.Range("AQ3").Value = 0.63016226204766679
$ws.Range("BM3").Value = 0.89719660353628139
$ws.Range("N4").Value = 0.87298131754730268
$ws.Range("AJ4").Value = 0.88071849153783943
$ws.Range("C5").Value = 0.57644898495068198
$ws.Range("D5").Value = 0.84686429042225653
$ws.Range("U5").Value = 0.90728848606264634
$ws.Range("AR5").Value = 0.87720831202463556
$ws.Range("AU6").Value = 0.95524159480490511
$ws.Range("G8").Value = 0.81327361820717869
$ws.Range("H9").Value = 0.90092936384962752
$ws.Range("AX9").Value = 0.82703502672572149
$ws.Range("BF9").Value = 0.91449960804514041
$ws.Range("H10").Value = 0.99795307822899404
$ws.Range("AD10").Value = 0.82993615350756056
$ws.Range("AQ10").Value = 0.90212141505601451
$ws.Range("M11").Value = 0.63132016125930834
$ws.Range("BE11").Value = 0.89904599781113081
$ws.Range("J12").Value = 0.99629009051279027
$ws.Range("AO12").Value = 0.94758918187317354
$ws.Range("BC12").Value = 0.7959248843981902
$ws.Range("AA13").Value = 0.71970318427108848
$ws.Range("B14").Value = 0.89249806798373355
$ws.Range("L14").Value = 0.78749697085718728
$ws.Range("AB14").Value = 0.67863715910227262
$ws.Range("AM14").Value = 0.9260536579214611
$ws.Range("L16").Value = 0.89704564295326827
$ws.Range("O16").Value = 0.89370813322193765
$ws.Range("G17").Value = 0.95871022822877272
$ws.Range("P17").Value = 0.95863796491937792
$ws.Range("D18").Value = 0.67855464335614579
$ws.Range("X18").Value = 0.96836932880891524
$ws.Range("AN18").Value = 0.97407421885286227
$ws.Range("AS18").Value = 0.69859279644369743
$ws.Range("I19").Value = 0.80350529744768795
$ws.Range("S20").Value = 0.8878425798241123
$ws.Range("V20").Value = 0.88946177601343246
$ws.Range("BN20").Value = 0.63565743680542508
$ws.Range("BJ21").Value = 0.94163315128629232
$ws.Range("Y22").Value = 0.73136606234353541
$ws.Range("BG22").Value = 0.74500527705930342
$ws.Range("AH23").Value = 0.98129822450643567
$ws.Range("AV23").Value = 0.94026971451611274
$ws.Range("AI24").Value = 0.98269528936771633
$ws.Range("BM25").Value = 0.88165350708234502
$ws.Range("B26").Value = 0.8782346350484791
$ws.Range("F26").Value = 0.93342335921566755
$ws.Range("G26").Value = 0.89492065022448808
$ws.Range("U26").Value = 0.6909581710248851
$ws.Range("W26").Value = 0.91526532334979882
$ws.Range("AV26").Value = 0.7206015693453991
$ws.Range("BI26").Value = 0.86942559235091399
$ws.Range("Y27").Value = 0.7446487859591755
$ws.Range("BN27").Value = 0.91798728310829047
$ws.Range("X28").Value = 0.87772655771169161
$ws.Range("AU28").Value = 0.71595148758294924
$ws.Range("BN29").Value = 0.69458089774907816
$ws.Range("AC30").Value = 0.94865234306784063
$ws.Range("BJ31").Value = 0.51853845942481613
$ws.Range("AO32").Value = 0.72809661572036422
$ws.Range("BI32").Value = 0.91160754073324579
$ws.Range("BH33").Value = 0.88670432313485947
$ws.Range("V35").Value = 0.81731617677670387
$ws.Range("AT35").Value = 0.96207420410719879
$ws.Range("AE36").Value = 0.94380545070175415
$ws.Range("AF36").Value = 0.96671542309323388
$ws.Range("BO36").Value = 0.77162660859959842
$ws.Range("F37").Value = 0.94167192198112248
$ws.Range("AN37").Value = 0.96082170198786532
$ws.Range("Z38").Value = 0.87505309610937532
$ws.Range("B39").Value = 0.73715423652978251
$ws.Range("AL39").Value = 0.9536563670096887
$ws.Range("BK39").Value = 0.7603149062887542
$ws.Range("AA40").Value = 0.84028929560522747
$ws.Range("AW40").Value = 0.70613048549901059
$ws.Range("AE41").Value = 0.96161684585017548
$ws.Range("AN41").Value = 0.97679288521595242
$ws.Range("AP43").Value = 0.91443564328676263
$ws.Range("AR43").Value = 0.92989951696268114
$ws.Range("BC43").Value = 0.7825088827082185
$ws.Range("I44").Value = 0.66505203510482491
$ws.Range("AA44").Value = 0.96260274082171449
$ws.Range("BA44").Value = 0.79619065855006066
$ws.Range("K45").Value = 0.92708424754956775
$ws.Range("AT45").Value = 0.56774773087332431
$ws.Range("T46").Value = 0.57177049947742964
$ws.Range("AA47").Value = 0.95136485687735806
$ws.Range("R48").Value = 0.93240090254374874
$ws.Range("AJ48").Value = 0.98783626072256081
$ws.Range("AC49").Value = 0.69702691332157074
$ws.Range("AQ49").Value = 0.82358294861521686
$ws.Range("AU49").Value = 0.8615414047331631
$ws.Range("BC50").Value = 0.71553762791031716
$ws.Range("F51").Value = 0.89204331440273266
$ws.Range("V51").Value = 0.8796523212102163
$ws.Range("BI51").Value = 0.82152272508798152
$ws.Range("BM51").Value = 0.8634733638998765
$ws.Range("AB52").Value = 0.851700771792403
$ws.Range("AH52").Value = 0.86416488341807118
$ws.Range("BE52").Value = 0.96921942877046607
$ws.Range("BB53").Value = 0.71804501505887708
$ws.Range("A54").Value = 0.72504575251778691
$ws.Range("AK54").Value = 0.78449624352049496
$ws.Range("AP54").Value = 0.98181460838763734
$ws.Range("BH54").Value = 0.97090032113643099
$ws.Range("Q55").Value = 0.97070896434629728
$ws.Range("K56").Value = 0.92845755334254365
$ws.Range("O56").Value = 0.93277668494328625
$ws.Range("BE56").Value = 0.85698283204900516
$ws.Range("AB57").Value = 0.76091431857363601
$ws.Range("E58").Value = 0.7656575772218448
$ws.Range("AV58").Value = 0.99492823097859406
$ws.Range("AY58").Value = 0.96805387968900369
$ws.Range("BO59").Value = 0.6286526321790411
$ws.Range("BI60").Value = 0.96715013126386162
$ws.Range("W61").Value = 0.76184606506844132
$ws.Range("BL62").Value = 0.7015893807254765
$ws.Range("J63").Value = 0.83547928803297211
$ws.Range("AD63").Value = 0.62796885496191535
$ws.Range("BL63").Value = 0.96620414162921842
$ws.Range("D64").Value = 0.93617951212426931
$ws.Range("AF65").Value = 0.669973025861059
$ws.Range("AQ65").Value = 0.9587744985679123
$ws.Range("BL65").Value = 0.93545373865221626
$ws.Range("AF66").Value = 0.8554643256101242
$ws.Range("BL66").Value = 0.96889324964846801
$ws.Range("BO66").Value = 0.94189224486309353
$ws.Range("AG67").Value = 0.90018862103075747
$ws.Range("A68").Value = 0.93154299352267222
$ws.Range("I68").Value = 0.9138445412754046
$ws.Range("X68").Value = 0.98593737175969609
$ws.Range("AI68").Value = 0.71385797640443549
